$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Locate the "User manager" bullet and the following bullet that
#    holds "i" + "habit_System " (two separate runs, no proofErr tags).
# ------------------------------------------------------------------
$searchRange = $d.Content
$found = $searchRange.Find.Execute("User manager", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$userManagerPara = $searchRange.Paragraphs(1)
$habitPara = $userManagerPara.Next()

$rangeStart = $userManagerPara.Range.Start
$rangeEnd = $habitPara.Range.End
$target = $d.Range($rangeStart, $rangeEnd)

# ------------------------------------------------------------------
# 2) Replace that two-paragraph span with:
#      - a single bullet paragraph where "ihabit_System" is wrapped in
#        proofErr spell-check markers (the trailing space moved to its
#        own run outside the markers)
#      - two new empty underline-formatted paragraphs
#    This mirrors the diff exactly (delete "User manager" paragraph,
#    rewrite the next one, and insert two blank lines after it).
# ------------------------------------------------------------------
$replacementXml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:r><w:t>habit_System</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr></w:p>
"@

[void]$target.InsertXML($replacementXml)

# ------------------------------------------------------------------
# 3) Remove the two blank underline paragraphs that used to sit right
#    after the "_GoBack" bookmark paragraph -- they were effectively
#    moved up next to the habit_System bullet in step 2.
# ------------------------------------------------------------------
$bookmark = $d.Bookmarks("_GoBack")
$bookmarkPara = $bookmark.Range.Paragraphs(1)

$firstExtra = $bookmarkPara.Next()
$firstExtra.Range.Delete()

$secondExtra = $bookmarkPara.Next()
$secondExtra.Range.Delete()
